$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy header formatting (bold font, border, alignment) from G1 to H1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Set the new header text
$ws.Range("H1").Value = "Save"

# Fill in the Save column values for rows 2-6
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0
